# Factura pdf y modificaciones de alquileres
# Rebuilds the "Pruebita" sheet's header + data rows for the new layout:
#   - merges "Fecha entrega temprano"/"Fecha entrega Tarde" into a single "Fecha entrega" column
#   - adds a new trailing "Precio" column
#   - updates rental rows 2-4 (new client/driver names, reserva, price) and
#     replaces rows 5-7 with new bookings (ids 4-6, new dates, prices, etc.)
# NumberFormat is forced to Text ("@") before writing any cell whose content looks
# like a plain number or a dd/mm/yyyy date, so Excel keeps storing it as text
# (matching the original workbook, where these values are plain strings too)
# instead of silently coercing it to a numeric/date cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Id"
$ws.Range("B1").Value = "Categoría"
$ws.Range("C1").Value = "Fecha recogida"
$ws.Range("D1").Value = "Ubicación recogida"
$ws.Range("E1").Value = "Ubicación entrega"
$ws.Range("F1").Value = "Fecha entrega"
$ws.Range("G1").Value = "Usuario cliente"
$ws.Range("H1").Value = "Contraseña cliente"
$ws.Range("I1").Value = "Conductores extra"
$ws.Range("J1").Value = "Usuario del conductor"
$ws.Range("K1").Value = "Contraseña del conductor"
$ws.Range("L1").Value = "Reserva"
$ws.Range("M1").Value = "Precio"

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "Pequeños"
$ws.Range("C2").Value = "15/11/23"
$ws.Range("D2").Value = "El camino"
$ws.Range("E2").Value = "Exostos"
$ws.Range("F2").Value = "16/11/23"
$ws.Range("G2").Value = "Tita"
$ws.Range("H2").Value = "Emp"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "Tita"
$ws.Range("K2").Value = "Loc"
$ws.Range("L2").Value = "Sí"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M2").Value = "1000"

$ws.Range("A3").NumberFormat = "@"
$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "SUV"
$ws.Range("C3").Value = "15/11/23"
$ws.Range("D3").Value = "Elegir"
$ws.Range("E3").Value = "Elegir"
$ws.Range("F3").Value = "16/11/23"
$ws.Range("G3").Value = "Martha"
$ws.Range("H3").Value = "Gen"
$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "1"
$ws.Range("J3").Value = "Martha"
$ws.Range("K3").Value = "Loc"
$ws.Range("L3").Value = "No"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M3").Value = "1000"

$ws.Range("A4").NumberFormat = "@"
$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "SUV"
$ws.Range("C4").Value = "15/11/23"
$ws.Range("D4").Value = "Ruedas"
$ws.Range("E4").Value = "Exostos"
$ws.Range("F4").Value = "16/11/23"
$ws.Range("G4").Value = "Mo"
$ws.Range("H4").Value = "Gen"
$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "1"
$ws.Range("J4").Value = "Mo"
$ws.Range("K4").Value = "Loc"
$ws.Range("L4").Value = "Sí"
$ws.Range("M4").NumberFormat = "@"
$ws.Range("M4").Value = "1000"

$ws.Range("A5").NumberFormat = "@"
$ws.Range("A5").Value = "4"
$ws.Range("B5").Value = "SUV"
$ws.Range("C5").NumberFormat = "@"
$ws.Range("C5").Value = "05/12/2023"
$ws.Range("D5").Value = "El camino"
$ws.Range("E5").Value = "Exostos"
$ws.Range("F5").NumberFormat = "@"
$ws.Range("F5").Value = "06/12/2023"
$ws.Range("G5").Value = "Juan"
$ws.Range("H5").Value = "Gen"
$ws.Range("I5").NumberFormat = "@"
$ws.Range("I5").Value = "1"
$ws.Range("J5").Value = "Juan"
$ws.Range("K5").Value = "Emp"
$ws.Range("L5").Value = "Sí"
$ws.Range("M5").NumberFormat = "@"
$ws.Range("M5").Value = "100"

$ws.Range("A6").NumberFormat = "@"
$ws.Range("A6").Value = "5"
$ws.Range("B6").Value = "SUV"
$ws.Range("C6").NumberFormat = "@"
$ws.Range("C6").Value = "05/12/2023"
$ws.Range("D6").Value = "El camino"
$ws.Range("E6").Value = "El camino"
$ws.Range("F6").NumberFormat = "@"
$ws.Range("F6").Value = "07/12/2023"
$ws.Range("G6").Value = "Juan"
$ws.Range("H6").Value = "Cli"
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value = "1"
$ws.Range("J6").Value = "Juan"
$ws.Range("K6").Value = "Emp"
$ws.Range("L6").Value = "Sí"
$ws.Range("M6").NumberFormat = "@"
$ws.Range("M6").Value = "700"

$ws.Range("A7").NumberFormat = "@"
$ws.Range("A7").Value = "6"
$ws.Range("B7").Value = "pequeños"
$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "10/12/2023"
$ws.Range("D7").Value = "Exostos"
$ws.Range("E7").Value = "ExostosXD"
$ws.Range("F7").NumberFormat = "@"
$ws.Range("F7").Value = "12/12/2023"
$ws.Range("G7").Value = "Juan"
$ws.Range("H7").Value = "Cli"
$ws.Range("I7").NumberFormat = "@"
$ws.Range("I7").Value = "0"
$ws.Range("J7").Value = "No"
$ws.Range("K7").Value = "NO"
$ws.Range("L7").Value = "Sí"
$ws.Range("M7").NumberFormat = "@"
$ws.Range("M7").Value = "12000"
